$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New regenerated s_val data (filtered save games), columns B:E and G for rows 2-9.
# Column order: TB, d2S, K, IP, Win(unchanged), sum
$data = @{
    2 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 10.19245300693656, 18.67282528286833)
    3 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 5.586269137925634)
    4 = @(0.2917716402565462, 0.002571899574220771, 261.3203778131603, 10.19245300693656, 271.8071743599276)
    5 = @(0.6606524410359556, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 6.348428708163715)
    6 = @(3.286832544864788, 1.655778082260271, 261.3203778131603, 10.19245300693656, 276.4554414472219)
    7 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 8.974608811992548)
    8 = @(3.286832544864788, 10.34677158129881, 3.537761648806719, 10.19245300693656, 27.36381878190688)
    9 = @(1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 3.754798637575387)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 7).Value = $vals[4]
}
